$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 1503.2258
$ws.Cells.Item(28, 10).Value = 1550.1111
$ws.Cells.Item(28, 12).Value = 1550.1111
$ws.Cells.Item(28, 14).Value = -2520.1111
$ws.Cells.Item(125, 8).Value = 2781.5557
$ws.Cells.Item(125, 9).Value = 1600
$ws.Cells.Item(125, 10).Value = 4258.5
$ws.Cells.Item(125, 11).Value = 14400
$ws.Cells.Item(125, 12).Value = 38326.5
$ws.Cells.Item(125, 13).Value = -11940
$ws.Cells.Item(125, 14).Value = -43246.5
$ws.Cells.Item(126, 8).Value = 158000
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 158000
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 158000
$ws.Cells.Item(126, 14).Value = -167880
$ws.Cells.Item(127, 8).Value = 886.25
$ws.Cells.Item(127, 9).Value = 886.25
$ws.Cells.Item(127, 10).Value = 0
$ws.Cells.Item(127, 11).Value = 2658.75
$ws.Cells.Item(127, 12).Value = 0
$ws.Cells.Item(127, 13).Value = 2301.25
$ws.Cells.Item(128, 8).Value = 0
$ws.Cells.Item(128, 9).Value = 0
$ws.Cells.Item(128, 10).Value = 0
$ws.Cells.Item(128, 11).Value = 0
$ws.Cells.Item(128, 12).Value = 0
$ws.Cells.Item(129, 8).Value = 3171.2222
$ws.Cells.Item(129, 9).Value = 1347
$ws.Cells.Item(129, 10).Value = 3692.4285
$ws.Cells.Item(129, 11).Value = 4041
$ws.Cells.Item(129, 12).Value = 11077.2855
$ws.Cells.Item(129, 13).Value = 959
$ws.Cells.Item(129, 14).Value = -21077.2855
$ws.Cells.Item(130, 8).Value = 89990
$ws.Cells.Item(130, 9).Value = 0
$ws.Cells.Item(130, 10).Value = 89990
$ws.Cells.Item(130, 11).Value = 0
$ws.Cells.Item(130, 12).Value = 89990
$ws.Cells.Item(130, 14).Value = -100030
$ws.Cells.Item(131, 8).Value = 4666.1113
$ws.Cells.Item(131, 9).Value = 1742.1428
$ws.Cells.Item(131, 10).Value = 14900
$ws.Cells.Item(131, 11).Value = 5226.428400000001
$ws.Cells.Item(131, 12).Value = 44700
$ws.Cells.Item(131, 13).Value = -186.4284000000007
$ws.Cells.Item(131, 14).Value = -54780
$ws.Cells.Item(132, 8).Value = 2583.678
$ws.Cells.Item(132, 9).Value = 2351.8235
$ws.Cells.Item(132, 10).Value = 4061.75
$ws.Cells.Item(132, 11).Value = 7055.470499999999
$ws.Cells.Item(132, 12).Value = 12185.25
$ws.Cells.Item(132, 13).Value = -4525.470499999999
$ws.Cells.Item(132, 14).Value = -17245.25
$ws.Cells.Item(133, 8).Value = 58333.168
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 58333.168
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 58333.168
$ws.Cells.Item(133, 14).Value = -68453.168
$ws.Cells.Item(134, 8).Value = 57939.53
$ws.Cells.Item(134, 9).Value = 0
$ws.Cells.Item(134, 10).Value = 57939.53
$ws.Cells.Item(134, 11).Value = 0
$ws.Cells.Item(134, 12).Value = 57939.53
$ws.Cells.Item(134, 14).Value = -68079.53
$ws.Cells.Item(135, 8).Value = 860.4074
$ws.Cells.Item(135, 9).Value = 778.11536
$ws.Cells.Item(135, 10).Value = 3000
$ws.Cells.Item(135, 11).Value = 7003.03824
$ws.Cells.Item(135, 12).Value = 27000
$ws.Cells.Item(135, 13).Value = -4468.03824
$ws.Cells.Item(135, 14).Value = -32070
$ws.Cells.Item(136, 8).Value = 45000
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 10).Value = 45000
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 12).Value = 45000
$ws.Cells.Item(136, 14).Value = -55200
$ws.Cells.Item(137, 8).Value = 3187.5
$ws.Cells.Item(137, 9).Value = 2087.15
$ws.Cells.Item(137, 10).Value = 5388.2
$ws.Cells.Item(137, 11).Value = 6261.450000000001
$ws.Cells.Item(137, 12).Value = 16164.6
$ws.Cells.Item(137, 13).Value = -3711.450000000001
$ws.Cells.Item(137, 14).Value = -21264.6
$ws.Cells.Item(138, 8).Value = 3923.8918
$ws.Cells.Item(138, 9).Value = 3721.4375
$ws.Cells.Item(138, 10).Value = 4078.1428
$ws.Cells.Item(138, 11).Value = 11164.3125
$ws.Cells.Item(138, 12).Value = 12234.4284
$ws.Cells.Item(138, 13).Value = -6024.3125
$ws.Cells.Item(138, 14).Value = -22514.4284
$ws.Cells.Item(139, 8).Value = 84523.22
$ws.Cells.Item(139, 9).Value = 200709
$ws.Cells.Item(139, 10).Value = 70000
$ws.Cells.Item(139, 11).Value = 200709
$ws.Cells.Item(139, 12).Value = 70000
$ws.Cells.Item(139, 13).Value = -195569
$ws.Cells.Item(139, 14).Value = -80280
$ws.Cells.Item(140, 8).Value = 99997
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 99997
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 99997
$ws.Cells.Item(140, 14).Value = -110357
$ws.Cells.Item(141, 8).Value = 5225.522
$ws.Cells.Item(141, 9).Value = 4188.6665
$ws.Cells.Item(141, 10).Value = 7169.625
$ws.Cells.Item(141, 11).Value = 12565.9995
$ws.Cells.Item(141, 12).Value = 21508.875
$ws.Cells.Item(141, 13).Value = -7385.999500000002
$ws.Cells.Item(141, 14).Value = -31868.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 3114
$ws.Cells.Item(122, 9).Value = 1285.2941
$ws.Cells.Item(122, 11).Value = 3855.8823
$ws.Cells.Item(122, 13).Value = -1405.8823

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2732.5
$ws.Cells.Item(31, 9).Value = 1980.7
$ws.Cells.Item(31, 10).Value = 3310.8076
$ws.Cells.Item(31, 11).Value = 1980.7
$ws.Cells.Item(31, 12).Value = 3310.8076
$ws.Cells.Item(31, 13).Value = -1685.7
$ws.Cells.Item(31, 14).Value = -3900.8076
$ws.Cells.Item(34, 8).Value = 2732.5
$ws.Cells.Item(34, 9).Value = 1980.7
$ws.Cells.Item(34, 10).Value = 3310.8076
$ws.Cells.Item(34, 11).Value = 1980.7
$ws.Cells.Item(34, 12).Value = 3310.8076
$ws.Cells.Item(34, 13).Value = -1778.7
$ws.Cells.Item(34, 14).Value = -3714.8076
$ws.Cells.Item(105, 8).Value = 1197.1666
$ws.Cells.Item(105, 9).Value = 1106.381
$ws.Cells.Item(105, 11).Value = 1106.381
$ws.Cells.Item(105, 13).Value = 640.6189999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(104, 8).Value = 75000
$ws.Cells.Item(104, 10).Value = 75000
$ws.Cells.Item(104, 12).Value = 75000
$ws.Cells.Item(104, 14).Value = -81988

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 1996.6487
$ws.Cells.Item(136, 9).Value = 1200.381
$ws.Cells.Item(136, 10).Value = 3041.75
$ws.Cells.Item(136, 11).Value = 3601.143
$ws.Cells.Item(136, 12).Value = 9125.25
$ws.Cells.Item(136, 13).Value = -1051.143
$ws.Cells.Item(136, 14).Value = -14225.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(119, 8).Value = 0
$ws.Cells.Item(119, 9).Value = 0
$ws.Cells.Item(119, 10).Value = 0
$ws.Cells.Item(119, 11).Value = 0
$ws.Cells.Item(119, 12).Value = 0
$ws.Cells.Item(120, 8).Value = 0
$ws.Cells.Item(120, 9).Value = 0
$ws.Cells.Item(120, 10).Value = 0
$ws.Cells.Item(120, 11).Value = 0
$ws.Cells.Item(120, 12).Value = 0
$ws.Cells.Item(121, 8).Value = 79999
$ws.Cells.Item(121, 9).Value = 0
$ws.Cells.Item(121, 10).Value = 79999
$ws.Cells.Item(121, 11).Value = 0
$ws.Cells.Item(121, 12).Value = 79999
$ws.Cells.Item(121, 14).Value = -83493
$ws.Cells.Item(122, 8).Value = 914990.6
$ws.Cells.Item(122, 9).Value = 1006089.8
$ws.Cells.Item(122, 10).Value = 3999
$ws.Cells.Item(122, 11).Value = 3018269.4
$ws.Cells.Item(122, 12).Value = 11997
$ws.Cells.Item(122, 13).Value = -3015819.4
$ws.Cells.Item(122, 14).Value = -16897
$ws.Cells.Item(123, 8).Value = 0
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 0
$ws.Cells.Item(123, 11).Value = 0
$ws.Cells.Item(123, 12).Value = 0
$ws.Cells.Item(124, 8).Value = 91583.336
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = 91583.336
$ws.Cells.Item(124, 11).Value = 0
$ws.Cells.Item(124, 12).Value = 91583.336
$ws.Cells.Item(124, 14).Value = -101403.336
$ws.Cells.Item(125, 8).Value = 237625
$ws.Cells.Item(125, 9).Value = 0
$ws.Cells.Item(125, 10).Value = 237625
$ws.Cells.Item(125, 11).Value = 0
$ws.Cells.Item(125, 12).Value = 237625
$ws.Cells.Item(125, 14).Value = -247465
$ws.Cells.Item(126, 8).Value = 2902
$ws.Cells.Item(126, 9).Value = 2288.7144
$ws.Cells.Item(126, 10).Value = 4333
$ws.Cells.Item(126, 11).Value = 6866.1432
$ws.Cells.Item(126, 12).Value = 12999
$ws.Cells.Item(126, 13).Value = -4396.1432
$ws.Cells.Item(126, 14).Value = -17939
$ws.Cells.Item(127, 8).Value = 64500
$ws.Cells.Item(127, 9).Value = 64500
$ws.Cells.Item(127, 10).Value = 0
$ws.Cells.Item(127, 11).Value = 64500
$ws.Cells.Item(127, 12).Value = 0
$ws.Cells.Item(127, 13).Value = -59540
$ws.Cells.Item(128, 8).Value = 0
$ws.Cells.Item(128, 9).Value = 0
$ws.Cells.Item(128, 10).Value = 0
$ws.Cells.Item(128, 11).Value = 0
$ws.Cells.Item(128, 12).Value = 0
$ws.Cells.Item(129, 8).Value = 99429
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 99429
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 12).Value = 99429
$ws.Cells.Item(129, 14).Value = -109429
$ws.Cells.Item(130, 8).Value = 0
$ws.Cells.Item(130, 9).Value = 0
$ws.Cells.Item(130, 10).Value = 0
$ws.Cells.Item(130, 11).Value = 0
$ws.Cells.Item(130, 12).Value = 0
$ws.Cells.Item(131, 8).Value = 42798.8
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = 42798.8
$ws.Cells.Item(131, 11).Value = 0
$ws.Cells.Item(131, 12).Value = 42798.8
$ws.Cells.Item(131, 14).Value = -52878.8
$ws.Cells.Item(132, 8).Value = 13333.728
$ws.Cells.Item(132, 9).Value = 13333.728
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 40001.18399999999
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -37471.18399999999
$ws.Cells.Item(133, 8).Value = 60357
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 60357
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 60357
$ws.Cells.Item(133, 14).Value = -70477
$ws.Cells.Item(135, 8).Value = 47690.92
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 47690.92
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 47690.92
$ws.Cells.Item(135, 14).Value = -57830.92
$ws.Cells.Item(136, 8).Value = 12370.75
$ws.Cells.Item(136, 9).Value = 12675.154
$ws.Cells.Item(136, 10).Value = 9996.4
$ws.Cells.Item(136, 11).Value = 38025.462
$ws.Cells.Item(136, 12).Value = 29989.2
$ws.Cells.Item(136, 13).Value = -35475.462
$ws.Cells.Item(136, 14).Value = -35089.2
$ws.Cells.Item(137, 8).Value = 66999.2
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 66999.2
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = 66999.2
$ws.Cells.Item(137, 14).Value = -77199.2
$ws.Cells.Item(138, 8).Value = 97017.664
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 97017.664
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 97017.664
$ws.Cells.Item(138, 14).Value = -107297.664
$ws.Cells.Item(139, 8).Value = 69937.5
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 69937.5
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 69937.5
$ws.Cells.Item(139, 14).Value = -80217.5
$ws.Cells.Item(140, 8).Value = 70500
$ws.Cells.Item(140, 9).Value = 70500
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 11).Value = 70500
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 13).Value = -65320
$ws.Cells.Item(141, 8).Value = 91998.78
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 91998.78
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 91998.78
$ws.Cells.Item(141, 14).Value = -102358.78
